# New submission synced: 2026-02-09 17:33:58
# Target sheet: "JSS 3A" (the sheet holding the Google-Form-style
# Timestamp / Full Name / Admission No / AI Score columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3A")

# The previous last row (row 8, Rachel Michael) had her Admission No
# typed in as text ("34"); normalize it to a real number to match the
# freshly-synced data below.
$ws.Cells.Item(8, 3).Value = 34

# Append the newly synced submission as row 9.
$ws.Cells.Item(9, 1).Value = "2026-02-09 17:33:58"
$ws.Cells.Item(9, 2).Value = "JOHANNAH SALEM BOURMANDA"

# Admission No for this submission stays as text ("28"), matching how
# the form sync originally wrote it before any normalization happened.
$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = "28"
$ws.Cells.Item(9, 3).Style = "Normal"

$ws.Cells.Item(9, 4).Value = 8
